# Updates league bases (Denmark Division 3) - re-shuffles several match rows
# whose B:AC data (id..PL_AhUnder) got reordered relative to the fixture
# list while column A (row index) stays put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowRange($row) {
    return $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 29))
}

# ---- Group 1: rows 32 <-> 34 (simple swap) ----
$r32 = Get-RowRange 32
$r34 = Get-RowRange 34
$v32 = $r32.Value2
$v34 = $r34.Value2
$r32.Value2 = $v34
$r34.Value2 = $v32

# ---- Group 2: rows 85,86,87,88,89 (5-cycle rotation) ----
# new_85 = old_89 ; new_86 = old_85 ; new_87 = old_86 ; new_88 = old_87 ; new_89 = old_88
$r85 = Get-RowRange 85
$r86 = Get-RowRange 86
$r87 = Get-RowRange 87
$r88 = Get-RowRange 88
$r89 = Get-RowRange 89
$v85 = $r85.Value2
$v86 = $r86.Value2
$v87 = $r87.Value2
$v88 = $r88.Value2
$v89 = $r89.Value2
$r85.Value2 = $v89
$r86.Value2 = $v85
$r87.Value2 = $v86
$r88.Value2 = $v87
$r89.Value2 = $v88

# ---- Group 3: rows 95 <-> 96 (simple swap) ----
$r95 = Get-RowRange 95
$r96 = Get-RowRange 96
$v95 = $r95.Value2
$v96 = $r96.Value2
$r95.Value2 = $v96
$r96.Value2 = $v95

# ---- Group 4: rows 131 <-> 132 (simple swap) ----
$r131 = Get-RowRange 131
$r132 = Get-RowRange 132
$v131 = $r131.Value2
$v132 = $r132.Value2
$r131.Value2 = $v132
$r132.Value2 = $v131

# ---- Group 5: rows 137 <-> 138 (simple swap) ----
$r137 = Get-RowRange 137
$r138 = Get-RowRange 138
$v137 = $r137.Value2
$v138 = $r138.Value2
$r137.Value2 = $v138
$r138.Value2 = $v137

# ---- Group 6: rows 139 <-> 140 (simple swap) ----
$r139 = Get-RowRange 139
$r140 = Get-RowRange 140
$v139 = $r139.Value2
$v140 = $r140.Value2
$r139.Value2 = $v140
$r140.Value2 = $v139

# ---- Group 7: rows 154 <-> 155 (simple swap) ----
$r154 = Get-RowRange 154
$r155 = Get-RowRange 155
$v154 = $r154.Value2
$v155 = $r155.Value2
$r154.Value2 = $v155
$r155.Value2 = $v154

# ---- Group 8: rows 170,171,172 (3-cycle rotation) ----
# new_170 = old_171 ; new_171 = old_172 ; new_172 = old_170
$r170 = Get-RowRange 170
$r171 = Get-RowRange 171
$r172 = Get-RowRange 172
$v170 = $r170.Value2
$v171 = $r171.Value2
$v172 = $r172.Value2
$r170.Value2 = $v171
$r171.Value2 = $v172
$r172.Value2 = $v170
